# Circle Language Spec Plan: Set font to Calibri for non-heading text.
#
# This script:
#   1. Changes the "Normal" paragraph style's font from Tahoma to Calibri
#      and gives it an explicit 11pt (22 half-points) size.
#   2. Normalizes the attribute order ("Month"/"Day"/"Year") inside the two
#      date smartTag's <w:smartTagPr> blocks in the "Date: ..." paragraph,
#      by rebuilding that paragraph's OOXML in place.

$d = $word.ActiveDocument

# --- 1. Normal style: Tahoma -> Calibri, explicit 11pt size -----------------
$normal = $d.Styles("Normal")
$normal.Font.Name = "Calibri"
$normal.Font.Size = 11

# --- 2. Reorder the smartTagPr date attributes (Month, Day, Year) ----------
# Locate the paragraph that holds the two "date" smart tags (the
# "Date: <date> - <date>" line) and rebuild it verbatim, except for the
# attribute order inside each <w:smartTagPr>.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "Date:*") {
        $datePara = $candidate
        break
    }
}

$newParaXml = @'
<w:p w:rsidR="002D5BF0" w:rsidRDefault="00D363DF"><w:pPr><w:jc w:val="center"/><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00772066"><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>Date</w:t></w:r><w:r w:rsidR="002D5BF0" w:rsidRPr="00772066"><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Month" w:val="6"/><w:attr w:name="Day" w:val="24"/><w:attr w:name="Year" w:val="2009"/></w:smartTagPr><w:r w:rsidR="00B56025"><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>June 24, 2009</w:t></w:r></w:smartTag><w:r w:rsidR="002D5BF0" w:rsidRPr="00772066"><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Month" w:val="9"/><w:attr w:name="Day" w:val="10"/><w:attr w:name="Year" w:val="2009"/></w:smartTagPr><w:r w:rsidR="003B53E5"><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>September 10, 2009</w:t></w:r></w:smartTag></w:p>
'@

if ($datePara) {
    $null = $datePara.Range.InsertXML($newParaXml)
}

Write-Output "done"
